$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the fighter/bomber building rows back to their "reverted" Chinese labels ---
# (reverting a merge that had renamed "fighter/bomber" buildings to "attack/escort" ones)
$ws.Range("B6").Value = "战斗机生产车间"
$ws.Range("B7").Value = "战斗机机库"
$ws.Range("B8").Value = "轰炸机生产车间"
$ws.Range("B9").Value = "轰炸机机库"

# --- Restore the view/selection state (active cell moves from B11 to A14) ---
[void]$ws.Range("A14").Select()
